# Automatic update of files.
# - Bump the "Förändrad" (changed) date in column C for all data rows (2-41) from 45705 to 45706.
# - Swap rows 39 and 40 (Beteckning + Area values).
# - Normalize row 41's height (matches the other data rows).
# - Append a new data row (42) for "A 7421-2025".
# - Grow the sheet dimension to include the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column C: 45705 -> 45706 for every data row (2 through 41) ---
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 3).Value = 45706
}

# --- 2) Swap the contents of row 39 and row 40 (A + G columns) ---
$a39 = $ws.Range("A39").Value()
$g39 = $ws.Range("G39").Value()
$a40 = $ws.Range("A40").Value()
$g40 = $ws.Range("G40").Value()

$ws.Range("A39").Value = $a40
$ws.Range("G39").Value = $g40
$ws.Range("A40").Value = $a39
$ws.Range("G40").Value = $g39

# --- 3) Row 41 picks up an explicit (default) row height, like the rows above it ---
$ws.Rows.Item(41).RowHeight = 15

# --- 4) Append new row 42 ---
$ws.Range("A42").Value = "A 7421-2025"
$ws.Range("B42").Value = 45702
$ws.Range("B42").NumberFormat = "YYYY-MM-DD"
$ws.Range("C42").Value = 45706
$ws.Range("C42").NumberFormat = "YYYY-MM-DD"
$ws.Range("D42").Value = "OKÄNT"
$ws.Range("E42").Value = "OKÄNT"
$ws.Range("F42").Value = "Kommuner"
$ws.Range("G42").Value = 2.9
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("R42").WrapText = $true
$ws.Range("R42").Value = ""
